$d = $word.ActiveDocument

# 1. Locate the word "Фамілія" (without the trailing ": ") so we know
#    exactly where it ends - the split point between the label and ": ".
$r = $d.Content
$r.Find.Execute("Фамілія", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# 2. Collapse the found range to its end and move the (singleton) _GoBack
#    bookmark there - this both relocates the bookmark away from the end
#    of the "Мова навчання" paragraph and forces a run split right after
#    "Фамілія", matching the target run layout.
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null

# 3. Rename the now-isolated "Фамілія" run to "Прізвище" (keeps the
#    Ukrainian-language run formatting; the surrounding ": " run is left
#    untouched).
$d.Content.Find.Execute("Фамілія", $true, $false, $false, $false, $false, $true, 1, $false, "Прізвище", 2) | Out-Null
